$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.833.85"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").Value = "2.292.59"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.513"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.78%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.66"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.36%  "

# Row 12
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.82"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.39%  "

# Row 15
$ws.Range("D15").Value = "2.645.23"
$ws.Range("E15").Value = "  -1.02%  "

# Row 16
$ws.Range("D16").Value = "2.289.16"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.99%  "

# Row 18
$ws.Range("D18").Value = "42.720.97"
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.21%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -0.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.43"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("E24").Value = "  -1.76%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.20%  "

# Row 27
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.93%  "

# Row 30
$ws.Range("E30").Value = "  -1.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.03"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.84"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.16%  "

# Row 33
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.87"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.81%  "

# Row 34
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.01"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.97"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.63%  "

# Row 37
$ws.Range("E37").Value = "  -1.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0685"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.40%  "

# Row 39
$ws.Range("E39").Value = "  -1.41%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.78%  "

# Row 43
$ws.Range("D43").Value = "2.018.80"
$ws.Range("E43").Value = "  +1.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.30%  "

# Row 49
$ws.Range("E49").Value = "  -2.36%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.20"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.83%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.514.22"
$ws.Range("E51").Value = "  -0.90%  "
